$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry the existing row 10 formatting (date style on column A, etc.) down
# into the new row 11 before filling in the new figures, same as a user
# would get by copying the prior row's formatting for consistency.
$ws.Range("A10:L10").Copy()
$ws.Range("A11:L11").PasteSpecial(-4122) | Out-Null

# Append row 11 with the new age-group case data (week of 02-Jan update)
$ws.Range("A11").Value = 44195
$ws.Range("B11").Value = 9814
$ws.Range("C11").Value = 9907
$ws.Range("D11").Value = 8822
$ws.Range("E11").Value = 7895
$ws.Range("F11").Value = 8695
$ws.Range("G11").Value = 6156
$ws.Range("H11").Value = 3121
$ws.Range("I11").Value = 2441
$ws.Range("J11").Value = 73
$ws.Range("K11").Value = 81
$ws.Range("L11").Value = 58.3

# Update the active selection like Excel would after entering the new row
$ws.Range("K20").Select()
